$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some Price (column D) values look like plain decimal numbers
# (e.g. "592.02") and Excel's COM layer would silently convert them to a
# floating point number, losing the original text formatting / precision.
# Force those specific cells to Text format before assigning so they stay
# as strings, exactly like the rest of column D (which already contains
# dotted "thousands" numbers that Excel never mis-detects, e.g. "63.140.36").

# --- Row 2 ---
$ws.Range("D2").Value = "63.140.36"
$ws.Range("E2").Value = "  -0.91%  "

# --- Row 3 ---
$ws.Range("D3").Value = "3.183.05"
$ws.Range("E3").Value = "  -3.95%  "

# --- Row 4 ---
$ws.Range("E4").Value = "  +0.04%  "

# --- Row 5 ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.02"
$ws.Range("E5").Value = "  -2.57%  "

# --- Row 6 ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.33"
$ws.Range("E6").Value = "  -4.51%  "

# --- Row 7 ---
$ws.Range("E7").Value = "  +0.04%  "

# --- Row 8 ---
$ws.Range("D8").Value = "3.179.45"
$ws.Range("E8").Value = "  -4.05%  "

# --- Row 9 ---
$ws.Range("E9").Value = "  -0.84%  "

# --- Row 10 ---
$ws.Range("E10").Value = "  -5.97%  "

# --- Row 11 ---
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.24"
$ws.Range("E11").Value = "  -5.42%  "

# --- Row 12 ---
$ws.Range("E12").Value = "  -3.22%  "

# --- Row 13 ---
$ws.Range("E13").Value = "  -4.38%  "

# --- Row 14 ---
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.67"
$ws.Range("E14").Value = "  -0.92%  "

# --- Row 15 ---
$ws.Range("D15").Value = "3.706.14"
$ws.Range("E15").Value = "  -3.92%  "

# --- Row 16 ---
$ws.Range("E16").Value = "  -1.07%  "

# --- Row 17 ---
$ws.Range("D17").Value = "3.185.59"
$ws.Range("E17").Value = "  -3.85%  "

# --- Row 18 ---
$ws.Range("D18").Value = "63.072.54"
$ws.Range("E18").Value = "  -1.11%  "

# --- Row 19 ---
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.57"
$ws.Range("E19").Value = "  -4.25%  "

# --- Row 20 ---
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.94"
$ws.Range("E20").Value = "  -4.00%  "

# --- Row 21 ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.05"
$ws.Range("E21").Value = "  -0.22%  "

# --- Row 22 ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("E22").Value = "  -5.64%  "

# --- Row 23 ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.64"
$ws.Range("E23").Value = "  -4.45%  "

# --- Row 24 ---
$ws.Range("E24").Value = "  -4.36%  "

# --- Row 25 ---
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.66"
$ws.Range("E25").Value = "  -2.95%  "

# --- Row 26 ---
$ws.Range("E26").Value = "  +0.00%  "

# --- Row 28 ---
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.67"
$ws.Range("E28").Value = "  -3.76%  "

# --- Row 29 & 30 swap: RenderToken <-> NEARProtocol ---
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.76"
$ws.Range("E29").Value = "  -5.67%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.67"
$ws.Range("E30").Value = "  -6.77%  "

# --- Row 31 ---
$ws.Range("E31").Value = "  -5.38%  "

# --- Row 32 ---
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.23"
$ws.Range("E32").Value = "  -5.68%  "

# --- Row 33 ---
$ws.Range("E33").Value = "  -3.21%  "

# --- Row 34 ---
$ws.Range("E34").Value = "  -5.65%  "

# --- Row 35 ---
$ws.Range("E35").Value = "  -5.92%  "

# --- Row 36 ---
$ws.Range("E36").Value = "  -4.18%  "

# --- Row 37 ---
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.35"
$ws.Range("E37").Value = "  -2.24%  "

# --- Row 38 ---
$ws.Range("E38").Value = "  -5.22%  "

# --- Row 39 ---
$ws.Range("E39").Value = "  -2.82%  "

# --- Row 40 ---
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "406.53"
$ws.Range("E40").Value = "  -6.38%  "

# --- Row 41 ---
$ws.Range("E41").Value = "  -3.10%  "

# --- Row 42 ---
$ws.Range("E42").Value = "  -3.30%  "

# --- Row 43 ---
$ws.Range("E43").Value = "  -5.89%  "

# --- Row 44 ---
$ws.Range("D44").Value = "2.813.40"
$ws.Range("E44").Value = "  -9.87%  "

# --- Row 45 ---
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.252"
$ws.Range("E45").Value = "  -5.66%  "

# --- Row 47 ---
$ws.Range("E47").Value = "  -5.32%  "

# --- Row 48 ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.17"
$ws.Range("E48").Value = "  -4.55%  "

# --- Row 49 ---
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.24"
$ws.Range("E49").Value = "  -0.34%  "

# --- Row 50 ---
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.24"
$ws.Range("E50").Value = "  -4.24%  "

# --- Row 51 ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.111"
$ws.Range("E51").Value = "  -1.99%  "
